# Natmi following Dr Hou advice
# Rebuild the LR-pair results table (header in row 1 stays untouched);
# rows 2-16 are replaced with the updated 5 (sender) x 3 (target) cluster matrix
# for the Adam10 -> Epha3 ligand-receptor pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 15,20
$arr[0,0] = "ECs"
$arr[0,1] = "Adam10"
$arr[0,2] = "Epha3"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 59.25485233333333
$arr[0,7] = 177.764557
$arr[0,8] = 0.2826627756245408
$arr[0,9] = 0.2826627756245408
$arr[0,10] = 2
$arr[0,11] = 0.6666666666666666
$arr[0,12] = 1.106217
$arr[0,13] = 3.318651
$arr[0,14] = 0.01813136626967656
$arr[0,15] = 0.01813136626967656
$arr[0,16] = 65.548724983623
$arr[0,17] = 589.938524852607
$arr[0,18] = 0.005125062315651954
$arr[0,19] = 0.005125062315651954
$arr[1,0] = "ECs"
$arr[1,1] = "Adam10"
$arr[1,2] = "Epha3"
$arr[1,3] = "FAPs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 59.25485233333333
$arr[1,7] = 177.764557
$arr[1,8] = 0.2826627756245408
$arr[1,9] = 0.2826627756245408
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 59.24481466666666
$arr[1,13] = 177.734444
$arr[1,14] = 0.971047664518299
$arr[1,15] = 0.971047664518299
$arr[1,16] = 3510.542744589034
$arr[1,17] = 31594.88470130131
$arr[1,18] = 0.2744790281164703
$arr[1,19] = 0.2744790281164703
$arr[2,0] = "ECs"
$arr[2,1] = "Adam10"
$arr[2,2] = "Epha3"
$arr[2,3] = "sCs"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 59.25485233333333
$arr[2,7] = 177.764557
$arr[2,8] = 0.2826627756245408
$arr[2,9] = 0.2826627756245408
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 0.6602006666666667
$arr[2,13] = 1.980602
$arr[2,14] = 0.01082096921202439
$arr[2,15] = 0.01082096921202439
$arr[2,16] = 39.12009301370156
$arr[2,17] = 352.080837123314
$arr[2,18] = 0.003058685192418513
$arr[2,19] = 0.003058685192418513
$arr[3,0] = "FAPs"
$arr[3,1] = "Adam10"
$arr[3,2] = "Epha3"
$arr[3,3] = "ECs"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 26.49803233333333
$arr[3,7] = 79.494097
$arr[3,8] = 0.1264032745503171
$arr[3,9] = 0.1264032745503171
$arr[3,10] = 2
$arr[3,11] = 0.6666666666666666
$arr[3,12] = 1.106217
$arr[3,13] = 3.318651
$arr[3,14] = 0.01813136626967656
$arr[3,15] = 0.01813136626967656
$arr[3,16] = 29.312573833683
$arr[3,17] = 263.813164503147
$arr[3,18] = 0.002291864068558284
$arr[3,19] = 0.002291864068558284
$arr[4,0] = "FAPs"
$arr[4,1] = "Adam10"
$arr[4,2] = "Epha3"
$arr[4,3] = "FAPs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 26.49803233333333
$arr[4,7] = 79.494097
$arr[4,8] = 0.1264032745503171
$arr[4,9] = 0.1264032745503171
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 59.24481466666666
$arr[4,13] = 177.734444
$arr[4,14] = 0.971047664518299
$arr[4,15] = 0.971047664518299
$arr[4,16] = 1569.871014619674
$arr[4,17] = 14128.83913157707
$arr[4,18] = 0.1227436045395507
$arr[4,19] = 0.1227436045395507
$arr[5,0] = "FAPs"
$arr[5,1] = "Adam10"
$arr[5,2] = "Epha3"
$arr[5,3] = "sCs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 26.49803233333333
$arr[5,7] = 79.494097
$arr[5,8] = 0.1264032745503171
$arr[5,9] = 0.1264032745503171
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 0.6602006666666667
$arr[5,13] = 1.980602
$arr[5,14] = 0.01082096921202439
$arr[5,15] = 0.01082096921202439
$arr[5,16] = 17.49401861182155
$arr[5,17] = 157.446167506394
$arr[5,18] = 0.001367805942208046
$arr[5,19] = 0.001367805942208046
$arr[6,0] = "M1"
$arr[6,1] = "Adam10"
$arr[6,2] = "Epha3"
$arr[6,3] = "ECs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 56.189477
$arr[6,7] = 168.568431
$arr[6,8] = 0.2680400491147057
$arr[6,9] = 0.2680400491147057
$arr[6,10] = 2
$arr[6,11] = 0.6666666666666666
$arr[6,12] = 1.106217
$arr[6,13] = 3.318651
$arr[6,14] = 0.01813136626967656
$arr[6,15] = 0.01813136626967656
$arr[6,16] = 62.15775467850901
$arr[6,17] = 559.419792106581
$arr[6,18] = 0.004859932305440824
$arr[6,19] = 0.004859932305440824
$arr[7,0] = "M1"
$arr[7,1] = "Adam10"
$arr[7,2] = "Epha3"
$arr[7,3] = "FAPs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 56.189477
$arr[7,7] = 168.568431
$arr[7,8] = 0.2680400491147057
$arr[7,9] = 0.2680400491147057
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 59.24481466666666
$arr[7,13] = 177.734444
$arr[7,14] = 0.971047664518299
$arr[7,15] = 0.971047664518299
$arr[7,16] = 3328.935151081929
$arr[7,17] = 29960.41635973736
$arr[7,18] = 0.2602796636902051
$arr[7,19] = 0.2602796636902051
$arr[8,0] = "M1"
$arr[8,1] = "Adam10"
$arr[8,2] = "Epha3"
$arr[8,3] = "sCs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 56.189477
$arr[8,7] = 168.568431
$arr[8,8] = 0.2680400491147057
$arr[8,9] = 0.2680400491147057
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 0.6602006666666667
$arr[8,13] = 1.980602
$arr[8,14] = 0.01082096921202439
$arr[8,15] = 0.01082096921202439
$arr[8,16] = 37.09633017505134
$arr[8,17] = 333.866971575462
$arr[8,18] = 0.002900453119059735
$arr[8,19] = 0.002900453119059735
$arr[9,0] = "M2"
$arr[9,1] = "Adam10"
$arr[9,2] = "Epha3"
$arr[9,3] = "ECs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 54.99878833333333
$arr[9,7] = 164.996365
$arr[9,8] = 0.2623601199583325
$arr[9,9] = 0.2623601199583325
$arr[9,10] = 2
$arr[9,11] = 0.6666666666666666
$arr[9,12] = 1.106217
$arr[9,13] = 3.318651
$arr[9,14] = 0.01813136626967656
$arr[9,15] = 0.01813136626967656
$arr[9,16] = 60.84059463373499
$arr[9,17] = 547.565351703615
$arr[9,18] = 0.004756947429520807
$arr[9,19] = 0.004756947429520807
$arr[10,0] = "M2"
$arr[10,1] = "Adam10"
$arr[10,2] = "Epha3"
$arr[10,3] = "FAPs"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 54.99878833333333
$arr[10,7] = 164.996365
$arr[10,8] = 0.2623601199583325
$arr[10,9] = 0.2623601199583325
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 59.24481466666666
$arr[10,13] = 177.734444
$arr[10,14] = 0.971047664518299
$arr[10,15] = 0.971047664518299
$arr[10,16] = 3258.393021699562
$arr[10,17] = 29325.53719529606
$arr[10,18] = 0.2547641817482796
$arr[10,19] = 0.2547641817482796
$arr[11,0] = "M2"
$arr[11,1] = "Adam10"
$arr[11,2] = "Epha3"
$arr[11,3] = "sCs"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 54.99878833333333
$arr[11,7] = 164.996365
$arr[11,8] = 0.2623601199583325
$arr[11,9] = 0.2623601199583325
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 0.6602006666666667
$arr[11,13] = 1.980602
$arr[11,14] = 0.01082096921202439
$arr[11,15] = 0.01082096921202439
$arr[11,16] = 36.31023672352556
$arr[11,17] = 326.79213051173
$arr[11,18] = 0.002838990780532141
$arr[11,19] = 0.002838990780532141
$arr[12,0] = "sCs"
$arr[12,1] = "Adam10"
$arr[12,2] = "Epha3"
$arr[12,3] = "ECs"
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 12.689751
$arr[12,7] = 38.069253
$arr[12,8] = 0.06053378075210392
$arr[12,9] = 0.06053378075210392
$arr[12,10] = 2
$arr[12,11] = 0.6666666666666666
$arr[12,12] = 1.106217
$arr[12,13] = 3.318651
$arr[12,14] = 0.01813136626967656
$arr[12,15] = 0.01813136626967656
$arr[12,16] = 14.037618281967
$arr[12,17] = 126.338564537703
$arr[12,18] = 0.001097560150504693
$arr[12,19] = 0.001097560150504693
$arr[13,0] = "sCs"
$arr[13,1] = "Adam10"
$arr[13,2] = "Epha3"
$arr[13,3] = "FAPs"
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 12.689751
$arr[13,7] = 38.069253
$arr[13,8] = 0.06053378075210392
$arr[13,9] = 0.06053378075210392
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 59.24481466666666
$arr[13,13] = 177.734444
$arr[13,14] = 0.971047664518299
$arr[13,15] = 0.971047664518299
$arr[13,16] = 751.8019461611481
$arr[13,17] = 6766.217515450332
$arr[13,18] = 0.05878118642379327
$arr[13,19] = 0.05878118642379327
$arr[14,0] = "sCs"
$arr[14,1] = "Adam10"
$arr[14,2] = "Epha3"
$arr[14,3] = "sCs"
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 12.689751
$arr[14,7] = 38.069253
$arr[14,8] = 0.06053378075210392
$arr[14,9] = 0.06053378075210392
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 0.6602006666666667
$arr[14,13] = 1.980602
$arr[14,14] = 0.01082096921202439
$arr[14,15] = 0.01082096921202439
$arr[14,16] = 8.377782070034
$arr[14,17] = 75.400038630306
$arr[14,18] = 0.0006550341778059508
$arr[14,19] = 0.0006550341778059508

$ws.Range("A2:T16").Value = $arr
